# Position encoder backward compatibility
#
# MB_POSITION_ENCODER_SCALING used to be documented further down the
# register map (register 217 / row 48). It's moved up to reuse the old,
# no-longer-relevant MB_VOLTAGE_TRIPS_DEPRECATED slot (register 207 /
# row 38), and its Notes column now records the backward-compatibility
# history for anyone relying on the old register number/name. The row
# it used to occupy is cleared back down to just its register number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 38 (register 207) only had "MB_VOLTAGE_TRIPS_DEPRECATED," in
# column B and nothing else. It becomes the live
# MB_POSITION_ENCODER_SCALING register definition.
$ws.Range("B38").Value = "MB_POSITION_ENCODER_SCALING,"
$ws.Range("C38").Value = "Tenths of a mm extension per position encoder count"
$ws.Range("D38").Value = "R/W"
$ws.Range("E38").Value = "0.1mm / count"
$ws.Range("G38").Value = "Set to zero to disable encoder. Formerly MB_VOLTAGE_TRIPS_DEPRECATED,"

# Row 48 (register 217) used to hold the Name/Description/R-W/Units for
# MB_POSITION_ENCODER_SCALING - that content now lives in row 38, so
# clear it back to just the register number in column A.
$ws.Range("B48:E48").ClearContents()

# Leave the cursor near the top of the edited register block.
$ws.Range("A39").Select()
